# Ajout draft mapping f595a2bd5e53be80aa00972cfd76eee4a5f7087b
#
# 1) Metadata sheet: bump the "Date" value to the new generation timestamp.
# 2) Elements sheet: add a new "Mapping: ..." column (AL) carrying the
#    business-mapping for the ROR DropZone extension, with a value only on
#    the Extension.value[x] row ("zonePoser").

$wb = $excel.ActiveWorkbook

# --- 1. Metadata!B8 = Date -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-12T09:15:29+00:00"

# --- 2. Elements!AL = new mapping column -----------------------------------
$ws = $wb.Worksheets.Item("Elements")

$lastCol = 38  # column AL (37 = AK was the previous last column)

# Header, row 1
$ws.Cells.Item(1, $lastCol).Value = "Mapping: Spécification métier vers l'extension ROR DropZone"
$ws.Cells.Item(1, $lastCol).Font.Bold = $true
$ws.Cells.Item(1, $lastCol).Interior.Color = $ws.Cells.Item(1, $lastCol - 1).Interior.Color
$ws.Cells.Item(1, $lastCol).VerticalAlignment = -4160
$ws.Cells.Item(1, $lastCol).WrapText = $true

# Data rows 2-5: no mapping recorded, mirror the blank "empty string" cells
# used throughout the rest of the sheet (e.g. column AK) rather than leaving
# a truly empty cell.
for ($r = 2; $r -le 5; $r++) {
    $ws.Cells.Item($r, $lastCol).Value = "'"
    $ws.Cells.Item($r, $lastCol).VerticalAlignment = -4160
    $ws.Cells.Item($r, $lastCol).WrapText = $true
}

# Data row 6 (Extension.value[x]) carries the actual mapping target.
$ws.Cells.Item(6, $lastCol).Value = "zonePoser"
$ws.Cells.Item(6, $lastCol).VerticalAlignment = -4160
$ws.Cells.Item(6, $lastCol).WrapText = $true

# Size the new column the same way the rest of the sheet was generated
# (auto-fit to its content).
$ws.Columns.Item($lastCol).ColumnWidth = 65.17578125
